$wb = $excel.ActiveWorkbook

# Sheet references
$wsDisease = $wb.Worksheets.Item("Disease_Synonymous")

# Append a new data row (A4:B4) to the Disease_Synonymous sheet
$wsDisease.Range("A4").Value = "BO_H72"
$wsDisease.Range("B4").Value = "H05"

# Make Disease_Synonymous the active/selected sheet (was HCAHeart_ID before),
# and move the selection to D10 on it.
$wsDisease.Select()
$wsDisease.Range("D10").Select()
